# Update auto scs (lamda_1, lamda_2), time in ms, and auto capacity values
# in the dict_compose_poisson_Cl2_proba8 table (Sheet1, rows 2-54, columns B-E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: @(RowNumber, B_lamda1, C_lamda2, D_key, E_prob)
$data = @(
    @(2, 33.94444444444444, 1.95, 0, 0.134),
    @(3, 33.94444444444444, 1.95, 2, 0.002),
    @(4, 33.94444444444444, 1.95, 3, 0.006),
    @(5, 33.94444444444444, 1.95, 4, 0.009000000000000001),
    @(6, 33.94444444444444, 1.95, 5, 0.013),
    @(7, 33.94444444444444, 1.95, 6, 0.04),
    @(8, 33.94444444444444, 1.95, 7, 0.044),
    @(9, 33.94444444444444, 1.95, 8, 0.062),
    @(10, 33.94444444444444, 1.95, 9, 0.054),
    @(11, 33.94444444444444, 1.95, 10, 0.038),
    @(12, 33.94444444444444, 1.95, 11, 0.029),
    @(13, 33.94444444444444, 1.95, 12, 0.03),
    @(14, 33.94444444444444, 1.95, 13, 0.02),
    @(15, 33.94444444444444, 1.95, 14, 0.03),
    @(16, 33.94444444444444, 1.95, 15, 0.034),
    @(17, 33.94444444444444, 1.95, 16, 0.029),
    @(18, 33.94444444444444, 1.95, 17, 0.032),
    @(19, 33.94444444444444, 1.95, 18, 0.032),
    @(20, 33.94444444444444, 1.95, 19, 0.035),
    @(21, 33.94444444444444, 1.95, 20, 0.029),
    @(22, 33.94444444444444, 1.95, 21, 0.024),
    @(23, 33.94444444444444, 1.95, 22, 0.025),
    @(24, 33.94444444444444, 1.95, 23, 0.022),
    @(25, 33.94444444444444, 1.95, 24, 0.013),
    @(26, 33.94444444444444, 1.95, 25, 0.021),
    @(27, 33.94444444444444, 1.95, 26, 0.023),
    @(28, 33.94444444444444, 1.95, 27, 0.02),
    @(29, 33.94444444444444, 1.95, 28, 0.01),
    @(30, 33.94444444444444, 1.95, 29, 0.014),
    @(31, 33.94444444444444, 1.95, 30, 0.01),
    @(32, 33.94444444444444, 1.95, 31, 0.015),
    @(33, 33.94444444444444, 1.95, 32, 0.01),
    @(34, 33.94444444444444, 1.95, 33, 0.009000000000000001),
    @(35, 33.94444444444444, 1.95, 34, 0.008),
    @(36, 33.94444444444444, 1.95, 35, 0.005),
    @(37, 33.94444444444444, 1.95, 36, 0.008),
    @(38, 33.94444444444444, 1.95, 37, 0.007),
    @(39, 33.94444444444444, 1.95, 38, 0.007),
    @(40, 33.94444444444444, 1.95, 39, 0.009000000000000001),
    @(41, 33.94444444444444, 1.95, 40, 0.003),
    @(42, 33.94444444444444, 1.95, 41, 0.006),
    @(43, 33.94444444444444, 1.95, 42, 0.003),
    @(44, 33.94444444444444, 1.95, 43, 0.003),
    @(45, 33.94444444444444, 1.95, 44, 0.005),
    @(46, 33.94444444444444, 1.95, 45, 0.004),
    @(47, 33.94444444444444, 1.95, 47, 0.002),
    @(48, 33.94444444444444, 1.95, 48, 0.002),
    @(49, 33.94444444444444, 1.95, 49, 0.004),
    @(50, 33.94444444444444, 1.95, 54, 0.001),
    @(51, 33.94444444444444, 1.95, 55, 0.001),
    @(52, 33.94444444444444, 1.95, 58, 0.001),
    @(53, 33.94444444444444, 1.95, 59, 0.001),
    @(54, 33.94444444444444, 1.95, 66, 0.001)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
